$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.893.12"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "1.724.22"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4767"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06127"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("D10").Value = "1.716.19"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06885"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5976"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.405"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "26.851.37"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9987"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "1.940.61"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.370"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.346"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.070"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.799"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.89%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.381"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.946"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07899"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.640"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.79%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9984"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.510"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.651"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01479"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3789"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.734"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1143"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05336"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.727"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.231"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
